$wb = $excel.ActiveWorkbook

# --- Sheet "Vendas Diárias": update the Dinheiro total and the Total Geral ---
$wsVendas = $wb.Worksheets.Item("Vendas Diárias")
$wsVendas.Range("B5").Value = 140
$wsVendas.Range("B9").Value = 10000657.5

# --- Sheet "Detalhes": append the two new "Venda Balcão" transactions ---
$wsDetalhes = $wb.Worksheets.Item("Detalhes")

# Move the footer row ("Sistema desenvolvido por ROBSON ALVES") down from
# row 12 to row 14 (leaving row 13 blank) so the new rows fit in between.
# Use a single-cell Copy so the style (s="3") travels with the value
# without touching the rest of the row.
$wsDetalhes.Range("A12").Copy($wsDetalhes.Range("A14"))
$wsDetalhes.Range("A12").Clear()

# Give the new value cells the same numeric style as the existing "valor"
# column (s="2") by copying a same-column cell before overwriting it.
$wsDetalhes.Range("B10").Copy($wsDetalhes.Range("B11"))
$wsDetalhes.Range("B10").Copy($wsDetalhes.Range("B12"))

# New row 11
$wsDetalhes.Range("A11").Value = "Venda Balcão"
$wsDetalhes.Range("B11").Value = 75
$wsDetalhes.Range("C11").Value = "29/08/2025 23:02:26"
$wsDetalhes.Range("D11").Value = "receita"
$wsDetalhes.Range("E11").Value = "Dinheiro"

# New row 12
$wsDetalhes.Range("A12").Value = "Venda Balcão"
$wsDetalhes.Range("B12").Value = 47
$wsDetalhes.Range("C12").Value = "29/08/2025 23:05:18"
$wsDetalhes.Range("D12").Value = "receita"
$wsDetalhes.Range("E12").Value = "Dinheiro"
